$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing header labels (spaces -> dots)
$ws.Range("B1").Value = "Current.Year.Salary"
$ws.Range("C1").Value = "Predicted.Salary"

# Add new header for the residual-as-percent-of-salary column
$ws.Range("E1").Value = "Percent.of.Salary.Off"

# Populate the new column with Residual / Current Year Salary for every data row
$ws.Range("E2:E434").Formula = "=D2/B2"

# Bake the formula results down to static values (matches the source data,
# which stores literal numbers rather than live formulas)
$ws.Range("E2:E434").Copy()
$ws.Range("E2:E434").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = 0
